# feat: add hash slot to DataEntity
#
# Inserts a new "data_checksum" column immediately before the "id" column
# on every sheet that shares the DataEntity-style schema
# (data_path, data_format, has_sample, has_reference, id, name, description):
#   - DataEntity
#   - AlignmentSet
#   - VariantSet
#   - MassSpectrometryResults
#   - Array

$wb = $excel.ActiveWorkbook

$targetSheets = @("DataEntity", "AlignmentSet", "VariantSet", "MassSpectrometryResults", "Array")

foreach ($sheetName in $targetSheets) {
    $ws = $wb.Worksheets.Item($sheetName)

    # "id" currently lives in column E on each of these sheets; insert a new
    # blank column there so id/name/description shift right by one (E->F,
    # F->G, G->H) and write the new header into the freed-up column.
    $ws.Range("E1").EntireColumn.Insert()
    $ws.Range("E1").Value = "data_checksum"
}
